# Horarios actualizados Línea 141 - 564
# Updates the three schedule sheets (LP1912, LP1912-215, 6203-6173) with the
# new scrape timestamp (06:53:44) and the additional rows it produced.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: LP1912
#   - header timestamp / row-count bump
#   - new row inserted right before the old last data row (old row 38 -> 39)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: 06:53:44"
$ws1.Range("A3").Value = "Total filas: 34"

$ws1.Rows.Item(38).Insert()
$ws1.Range("A38").Value = "06:53:44"
$ws1.Range("B38").Value = "08:43"
$ws1.Range("C38").Value = "215C_EL PATO"
$ws1.Range("D38").Value = 110
$ws1.Range("E38").Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215
#   - header timestamp / row-count bump
#   - new row inserted right before the old last data row (old row 13 -> 14)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 06:53:44"
$ws2.Range("A3").Value = "Total filas: 9"

$ws2.Rows.Item(13).Insert()
$ws2.Range("A13").Value = "06:53:44"
$ws2.Range("B13").Value = "08:43"
$ws2.Range("C13").Value = "215C_EL PATO"
$ws2.Range("D13").Value = 110
$ws2.Range("E13").Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173
#   - header timestamp / row-count bump
#   - new row inserted right before the old last data row (old row 8 -> 9)
#   - brand-new row appended afterwards (row 10)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 06:53:44"
$ws3.Range("A3").Value = "Total filas: 5"

$ws3.Rows.Item(8).Insert()
$ws3.Range("A8").Value = "06:53:44"
$ws3.Range("B8").Value = "08:35"
$ws3.Range("C8").Value = "215A_LA PLATA"
$ws3.Range("D8").Value = 102
$ws3.Range("E8").Value = "L6173"

$ws3.Range("A10").Value = "06:53:44"
$ws3.Range("B10").Value = "08:50"
$ws3.Range("C10").Value = "215C_LA PLATA"
$ws3.Range("D10").Value = 117
$ws3.Range("E10").Value = "L6203"
